$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.168195009231567
$ws.Range("B1").Value = 2.278706550598145
$ws.Range("C1").Value = 2.889773607254028
$ws.Range("D1").Value = 2.134241104125977
$ws.Range("E1").Value = 2.067207336425781
